$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell C10 ("R30" rule's "From" value) changes from 18 to 1.
$ws.Range("C10").Value = 1
